# httk-benchmarks.xlsx: add invivoPKfit outputs row (v2.3.0) to the Sheet1 table,
# matching author commit "Added invivoPKfit outputs to dashboard script".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow Table1 by one row (A1:R23 -> A1:R24); ListRows.Add() appends a blank
# row right after the table's current last row and extends ref/autoFilter.
$lo = $ws.ListObjects.Item("Table1")
$newRow = $lo.ListRows.Add()

# New row is row 24 (table had header row 1 + 22 data rows, now 23 data rows).
$rowNum = $newRow.Range.Row

# Match the formatting used by the rest of the data rows (left-aligned).
$ws.Range("A" + $rowNum + ":R" + $rowNum).HorizontalAlignment = -4131

$ws.Range("A" + $rowNum).Value = "2.3.0"
$ws.Range("B" + $rowNum).Value = 1023
$ws.Range("C" + $rowNum).Value = 0.9999
$ws.Range("D" + $rowNum).Value = 1
$ws.Range("E" + $rowNum).Value = 1
$ws.Range("F" + $rowNum).Value = 1.063
$ws.Range("G" + $rowNum).Value = 352
$ws.Range("H" + $rowNum).Value = 0.2996
$ws.Range("I" + $rowNum).Value = 352
$ws.Range("J" + $rowNum).Value = 1.419
$ws.Range("K" + $rowNum).Value = 86
$ws.Range("L" + $rowNum).Value = 1.047
$ws.Range("M" + $rowNum).Value = 86
$ws.Range("N" + $rowNum).Value = 1.33
$ws.Range("O" + $rowNum).Value = 86
$ws.Range("P" + $rowNum).Value = 0.6344
$ws.Range("Q" + $rowNum).Value = 863
$ws.Range("R" + $rowNum).Value = "Used Caco-2 to replace Fabs=Fgut=1"

# Move the view down to the newly added row, like the author's saved file.
$ws.Range("R" + $rowNum).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 4
